# Additional options for filter overriding.
#
# - Mark the "check export to Excel ..." to-do item (row 27) as complete
#   by switching its row style from Neutral (pending) to Good (done).
# - Insert a new pending to-do row directly below it for:
#   "test each of the new data type examples with visual totals."
#   (this naturally shifts every following row down by one, which also
#   matches the updated sheet dimension of B1:H36).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mark the existing task as done.
$ws.Range("B27:F27").Style = "Good"

# Insert a new row right after it; Excel inherits the row-above's
# formatting (Neutral) for the newly inserted cells.
$ws.Rows("28:28").Insert()

# Make sure the new row has the expected style and fill in its text.
$ws.Range("B28:F28").Style = "Neutral"
$ws.Range("B28").Value = "test each of the new data type examples with visual totals."
